$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Columns D (Price) and E (Volume(1h)) are stored as text in the source data
# (e.g. "62.561.07", "  -3.49%  "), so NumberFormat is forced to Text ("@")
# before assignment to prevent Excel from auto-converting numeric-looking
# strings (like "3.92" or "0.0000149") into real numbers.


# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.561.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.49%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.011.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.98%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.95%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.58%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.569'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.43%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.014.94'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.87%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.28'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.32%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.367'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.535.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.23%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.50%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.661.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.03%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.76'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.13%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.009.84'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.15%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000149'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.59%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '394.58'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.91%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.69'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.48%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.25%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.78%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.469'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.41%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.93%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0967'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.60%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.51'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.27%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.54'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.47%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.89'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.35%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.68'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.98%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.02'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.79%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.08'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.72%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.29'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.45%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.57'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.71%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.454.68'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -10.07%  '

# Row 40
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.67'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.32%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.92'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.65%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.49'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.78%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.664'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.46%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0596'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.44%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.10%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0247'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.26%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -10.76%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0954'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.69%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.85'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.43%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.50'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '264.70'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -8.32%  '
